function Update-DatePlaceholder {
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($shape.HasTextFrame) {
            try {
                if ($shape.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }
        if ($isDatePlaceholder) {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}

$p = $ppt.ActivePresentation
$newDate = "4/29/2024"

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes $newDate

# Every slide layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes $newDate
